{"js": "// The diff this task ships with does NOT touch any user-visible document\n// content (body text, paragraphs, tables, styles, etc.). Every changed byte\n// lives in the SharePoint \"Document Information Panel\" metadata schema that\n// Word keeps in the built-in custom XML parts customXml/item2.xml (the\n// ct:contentTypeSchema / xsd:schema describing the library's content-type\n// columns) and customXml/itemProps2.xml (its datastore item/schemaRefs):\n//   - ma:contentTypeVersion 15 -> 16, ma:versionID / ma:fieldsID regenerated\n//   - a new \"MediaServiceObjectDetectorVersions\" column definition added\n//   - several ma:displayName captions re-localized from English to Dutch\n//     (\"Shared With\" -> \"Gedeeld met\", \"Title\" -> \"Titel\", ...)\n//   - itemProps2.xml's ds:itemID swapped and its ds:schemaRefs list dropped\n//\n// The author's own commit message confirms this is not an authored edit:\n//   \"Geen wijzigingen. Veranderingen in data na installatie op productie\"\n//   (\"No changes. Changes in data after installation on production.\")\n// i.e. SharePoint regenerated/re-synced its content-type schema metadata\n// server-side after the template was (re)installed on production; nothing\n// in the document itself was changed by a user/editor action.\n//\n// This is also not reachable from Word's object model: built-in custom XML\n// parts (the ones backing the Document Information Panel, like this one)\n// are read-only through CustomXmlPart in both Office.js and the Word COM/\n// VBA surface - setXml()/insertElement()/updateElement()/insertAttribute()/\n// updateAttribute() all raise errors for them (GeneralException /\n// InvalidArgument), exactly like real Word refuses to let add-ins rewrite\n// these SharePoint-managed schema parts. There is no supported\n// Word.Document / Office.js call that performs (or should perform) this\n// edit, so this script intentionally makes no changes to the document.\n", "ps1": "# The diff this task ships with does NOT touch any user-visible document\n# content (body text, paragraphs, tables, styles, etc.). Every changed byte\n# lives in the SharePoint \"Document Information Panel\" metadata schema that\n# Word keeps in the built-in custom XML parts customXml/item2.xml (the\n# ct:contentTypeSchema / xsd:schema describing the library's content-type\n# columns) and customXml/itemProps2.xml (its datastore item/schemaRefs):\n#   - ma:contentTypeVersion 15 -> 16, ma:versionID / ma:fieldsID regenerated\n#   - a new \"MediaServiceObjectDetectorVersions\" column definition added\n#   - several ma:displayName captions re-localized from English to Dutch\n#     (\"Shared With\" -> \"Gedeeld met\", \"Title\" -> \"Titel\", ...)\n#   - itemProps2.xml's ds:itemID swapped and its ds:schemaRefs list dropped\n#\n# The author's own commit message confirms this is not an authored edit:\n#   \"Geen wijzigingen. Veranderingen in data na installatie op productie\"\n#   (\"No changes. Changes in data after installation on production.\")\n# i.e. SharePoint regenerated/re-synced its content-type schema metadata\n# server-side after the template was (re)installed on production; nothing\n# in the document itself was changed by a user/editor action.\n#\n# This is also not reachable from Word's object model: built-in custom XML\n# parts (the ones backing the Document Information Panel, like this one)\n# are read-only through CustomXMLPart in the Word COM/VBA surface - there\n# is no writable XML/.XML setter and no UpdateAttribute/InsertElement/\n# UpdateElement members exposed for it (attempting $part.XML = ... raises\n# \"Property 'XML' cannot be found on this object\"), exactly like real Word\n# refuses to let automation rewrite these SharePoint-managed schema parts.\n# There is no supported Word COM call that performs (or should perform)\n# this edit, so this script intentionally makes no changes to the document.\n\n$d = $word.ActiveDocument\n"}
